$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 7).Value = 4.2
$ws.Cells.Item(2, 9).Value = 1.7
$ws.Cells.Item(2, 10).Value = 4.65
$ws.Cells.Item(2, 19).Value = 1.39
$ws.Cells.Item(2, 20).Value = 2.45
$ws.Cells.Item(2, 24).Value = 18.5
$ws.Cells.Item(2, 31).Value = 13
$ws.Cells.Item(2, 33).Value = 400
# Row 3
$ws.Cells.Item(3, 7).Value = 2.22
$ws.Cells.Item(3, 9).Value = 3
$ws.Cells.Item(3, 15).Value = 1.88
$ws.Cells.Item(3, 23).Value = 6.7
$ws.Cells.Item(3, 24).Value = 9.5
$ws.Cells.Item(3, 26).Value = 18.5
$ws.Cells.Item(3, 31).Value = 10.75
$ws.Cells.Item(3, 34).Value = 7.7
# Row 5
$ws.Cells.Item(5, 7).Value = 3.3
$ws.Cells.Item(5, 8).Value = 3.5
$ws.Cells.Item(5, 9).Value = 2.05
$ws.Cells.Item(5, 10).Value = 3.75
$ws.Cells.Item(5, 11).Value = 2.3
$ws.Cells.Item(5, 12).Value = 2.63
$ws.Cells.Item(5, 13).Value = 1.2
$ws.Cells.Item(5, 14).Value = 4.33
$ws.Cells.Item(5, 15).Value = 1.7
$ws.Cells.Item(5, 16).Value = 2.1
$ws.Cells.Item(5, 17).Value = 2.63
$ws.Cells.Item(5, 18).Value = 1.44
$ws.Cells.Item(5, 20).Value = 3.25
$ws.Cells.Item(5, 21).Value = 1.62
$ws.Cells.Item(5, 22).Value = 2.2
$ws.Cells.Item(5, 23).Value = 13
$ws.Cells.Item(5, 24).Value = 19
$ws.Cells.Item(5, 29).Value = 13
$ws.Cells.Item(5, 30).Value = 7
$ws.Cells.Item(5, 31).Value = 12
$ws.Cells.Item(5, 34).Value = 9.5
$ws.Cells.Item(5, 37).Value = 19
$ws.Cells.Item(5, 38).Value = 15
$ws.Cells.Item(5, 39).Value = 23
$ws.Cells.Item(5, 40).Value = 1.04
$ws.Cells.Item(5, 41).Value = 13
# Row 6
$ws.Cells.Item(6, 7).Value = 2.3
$ws.Cells.Item(6, 8).Value = 3.1
$ws.Cells.Item(6, 10).Value = 2.8
$ws.Cells.Item(6, 11).Value = 2.1
$ws.Cells.Item(6, 12).Value = 3.5
$ws.Cells.Item(6, 13).Value = 1.32
$ws.Cells.Item(6, 14).Value = 2.85
$ws.Cells.Item(6, 15).Value = 1.93
$ws.Cells.Item(6, 16).Value = 1.7
$ws.Cells.Item(6, 17).Value = 3.1
$ws.Cells.Item(6, 18).Value = 1.27
$ws.Cells.Item(6, 22).Value = 1.9
$ws.Cells.Item(6, 23).Value = 7.7
$ws.Cells.Item(6, 24).Value = 11.25
$ws.Cells.Item(6, 29).Value = 9
$ws.Cells.Item(6, 30).Value = 6.1
$ws.Cells.Item(6, 34).Value = 8.75
$ws.Cells.Item(6, 38).Value = 27
$ws.Cells.Item(6, 39).Value = 35
# Row 7
$ws.Cells.Item(7, 7).Value = 3.1
$ws.Cells.Item(7, 8).Value = 3.4
$ws.Cells.Item(7, 9).Value = 2.12
$ws.Cells.Item(7, 10).Value = 3.55
$ws.Cells.Item(7, 12).Value = 2.7
$ws.Cells.Item(7, 13).Value = 1.26
$ws.Cells.Item(7, 15).Value = 1.78
$ws.Cells.Item(7, 16).Value = 1.83
$ws.Cells.Item(7, 17).Value = 2.8
$ws.Cells.Item(7, 18).Value = 1.33
$ws.Cells.Item(7, 23).Value = 10.5
$ws.Cells.Item(7, 24).Value = 17
$ws.Cells.Item(7, 25).Value = 11
$ws.Cells.Item(7, 26).Value = 40
$ws.Cells.Item(7, 27).Value = 25
$ws.Cells.Item(7, 28).Value = 32
$ws.Cells.Item(7, 29).Value = 11
$ws.Cells.Item(7, 30).Value = 6.7
$ws.Cells.Item(7, 31).Value = 13.5
$ws.Cells.Item(7, 32).Value = 55
$ws.Cells.Item(7, 33).Value = 400
$ws.Cells.Item(7, 35).Value = 10.5
$ws.Cells.Item(7, 36).Value = 8.75
$ws.Cells.Item(7, 37).Value = 20
$ws.Cells.Item(7, 38).Value = 16.5
$ws.Cells.Item(7, 39).Value = 26
# Row 8
$ws.Cells.Item(8, 7).Value = 2.02
$ws.Cells.Item(8, 8).Value = 3.4
$ws.Cells.Item(8, 10).Value = 2.62
$ws.Cells.Item(8, 11).Value = 2.1
$ws.Cells.Item(8, 13).Value = 1.28
$ws.Cells.Item(8, 15).Value = 1.82
$ws.Cells.Item(8, 19).Value = 1.37
$ws.Cells.Item(8, 20).Value = 2.5
$ws.Cells.Item(8, 21).Value = 1.74
$ws.Cells.Item(8, 22).Value = 1.98
$ws.Cells.Item(8, 23).Value = 6.6
$ws.Cells.Item(8, 26).Value = 15
$ws.Cells.Item(8, 27).Value = 13
$ws.Cells.Item(8, 29).Value = 10
$ws.Cells.Item(8, 30).Value = 5.8
$ws.Cells.Item(8, 31).Value = 11.75
$ws.Cells.Item(8, 32).Value = 50
$ws.Cells.Item(8, 34).Value = 8.25
$ws.Cells.Item(8, 35).Value = 13.5
$ws.Cells.Item(8, 36).Value = 9.5
$ws.Cells.Item(8, 37).Value = 30
$ws.Cells.Item(8, 38).Value = 22
$ws.Cells.Item(8, 39).Value = 27
# Row 9
$ws.Cells.Item(9, 7).Value = 1.11
$ws.Cells.Item(9, 9).Value = 17
$ws.Cells.Item(9, 10).Value = 1.36
$ws.Cells.Item(9, 11).Value = 3.25
$ws.Cells.Item(9, 12).Value = 11.5
$ws.Cells.Item(9, 24).Value = 6.6
$ws.Cells.Item(9, 25).Value = 9.5
$ws.Cells.Item(9, 26).Value = 6.2
$ws.Cells.Item(9, 29).Value = 23
$ws.Cells.Item(9, 32).Value = 80
$ws.Cells.Item(9, 33).Value = 450
$ws.Cells.Item(9, 34).Value = 55
$ws.Cells.Item(9, 35).Value = 175
$ws.Cells.Item(9, 36).Value = 50
$ws.Cells.Item(9, 38).Value = 200
$ws.Cells.Item(9, 39).Value = 110
# Row 11
$ws.Cells.Item(11, 40).Value = 1.03
$ws.Cells.Item(11, 41).Value = 10
# Row 12
$ws.Cells.Item(12, 7).Value = 2.35
$ws.Cells.Item(12, 9).Value = 2.7
$ws.Cells.Item(12, 24).Value = 13
$ws.Cells.Item(12, 27).Value = 17
$ws.Cells.Item(12, 37).Value = 29
$ws.Cells.Item(12, 38).Value = 21
# Row 13
$ws.Cells.Item(13, 13).Value = 1.3
$ws.Cells.Item(13, 14).Value = 3.4
$ws.Cells.Item(13, 15).Value = 1.98
$ws.Cells.Item(13, 16).Value = 1.83
$ws.Cells.Item(13, 17).Value = 3.4
$ws.Cells.Item(13, 18).Value = 1.3
$ws.Cells.Item(13, 33).Value = 600
$ws.Cells.Item(13, 40).Value = 1.06
$ws.Cells.Item(13, 41).Value = 8
# Row 14
$ws.Cells.Item(14, 7).Value = 1.73
$ws.Cells.Item(14, 10).Value = 2.4
$ws.Cells.Item(14, 17).Value = 4.5
$ws.Cells.Item(14, 18).Value = 1.18
$ws.Cells.Item(14, 30).Value = 7
$ws.Cells.Item(14, 34).Value = 11
$ws.Cells.Item(14, 44).Value = 1.78
$ws.Cells.Item(14, 45).Value = 2.03
# Row 18
$ws.Cells.Item(18, 10).Value = 1.88
$ws.Cells.Item(18, 12).Value = 7.9
$ws.Cells.Item(18, 23).Value = 5.1
# Row 20
$ws.Cells.Item(20, 7).Value = 7.2
$ws.Cells.Item(20, 9).Value = 1.45
$ws.Cells.Item(20, 10).Value = 7.1
$ws.Cells.Item(20, 12).Value = 2.05
$ws.Cells.Item(20, 21).Value = 2.42
$ws.Cells.Item(20, 22).Value = 1.44
$ws.Cells.Item(20, 25).Value = 26
$ws.Cells.Item(20, 30).Value = 8
$ws.Cells.Item(20, 34).Value = 4.75
$ws.Cells.Item(20, 35).Value = 5.4
